$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 174; this shifts the former rows 174:203 down to 175:204
$ws.Rows.Item(174).Insert()

# Populate the new row 174 with a fresh weekly record (same fixed attributes as
# the record that used to occupy row 174, with updated date/volume/price figures)
$ws.Range("A174").Value = 1
$ws.Range("B174").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C174").Value = "Arica y Parinacota"
$ws.Range("D174").Value = 44491
$ws.Range("E174").Value = 15
$ws.Range("F174").Value = 100114013
$ws.Range("G174").Value = "Zanahoria"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 100
$ws.Range("K174").Value = 12000
$ws.Range("L174").Value = 13000
$ws.Range("M174").Value = 12500
$ws.Range("N174").Value = "$/saco 25 kilos"
$ws.Range("O174").Value = "Región de Arica y Parinacota"
$ws.Range("P174").Value = 500
$ws.Range("Q174").Value = 25
$ws.Range("R174").Value = "Hortaliza"

# Apply the same date format style used by the other rows' Fecha column
$ws.Range("D174").NumberFormat = $ws.Range("D175").NumberFormat
